$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell L1, styled like the other header cells (bold)
$ws.Range("L1").Value = "Note:"
$ws.Range("L1").Font.Bold = $true

# New Trawl numbers in column A
$ws.Range("A2").Value = 15
$ws.Range("A3").Value = 14

# New "Used as voucher" notes in column L
$ws.Range("L2").Value = "Used as voucher"
$ws.Range("L3").Value = "Used as voucher"

# Move the active selection, matching the saved cursor position in the diff
$ws.Range("L4").Select()
